# This script applies the recomputed NATMI TPM values for the Clcf1 -> Cntfr
# ligand-receptor pair sheet (rows 2-10, columns E-T), matching the updated
# TPM-based recalculation referenced in the commit "update scripts wuth new tpm".
# Values are applied as literal numbers (the sheet has no formulas).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> hashtable of column letter -> new numeric value (as string,
# parsed with [double] to safely handle scientific notation like "E-05").
$rowUpdates = @{}

$rowUpdates[2] = @{
    "E" = "1";
    "F" = "0.3333333333333333";
    "G" = "0.06373566666666666";
    "H" = "0.191207";
    "I" = "0.01058875298517695";
    "J" = "0.01058875298517695";
    "M" = "0.04154133333333333";
    "N" = "0.124624";
    "O" = "0.002517093804502335";
    "P" = "0.002517093804502335";
    "Q" = "0.002647664574222222";
    "R" = "0.023828981168";
    "S" = "2.665288453639451E-05";
    "T" = "2.665288453639451E-05"
}

$rowUpdates[3] = @{
    "E" = "1";
    "F" = "0.3333333333333333";
    "G" = "0.06373566666666666";
    "H" = "0.191207";
    "I" = "0.01058875298517695";
    "J" = "0.01058875298517695";
    "O" = "0.8133013372545576";
    "P" = "0.8133013372545578";
    "Q" = "0.8554902224798888";
    "R" = "7.699412002319";
    "S" = "0.008611846962702605";
    "T" = "0.008611846962702606"
}

$rowUpdates[4] = @{
    "E" = "1";
    "F" = "0.3333333333333333";
    "G" = "0.06373566666666666";
    "H" = "0.191207";
    "I" = "0.01058875298517695";
    "J" = "0.01058875298517695";
    "O" = "0.18418156894094";
    "P" = "0.18418156894094";
    "Q" = "0.1937357338202222";
    "R" = "1.743621604382";
    "S" = "0.001950253137937953";
    "T" = "0.001950253137937953"
}

$rowUpdates[5] = @{
    "I" = "0.2961697031425515";
    "J" = "0.2961697031425515";
    "M" = "0.04154133333333333";
    "N" = "0.124624";
    "O" = "0.002517093804502335";
    "P" = "0.002517093804502335";
    "Q" = "0.07405574878044445";
    "R" = "0.666501739024";
    "S" = "0.0007454869248614121";
    "T" = "0.0007454869248614121"
}

$rowUpdates[6] = @{
    "I" = "0.2961697031425515";
    "J" = "0.2961697031425515";
    "O" = "0.8133013372545576";
    "P" = "0.8133013372545578";
    "S" = "0.2408752156201225";
    "T" = "0.2408752156201226"
}

$rowUpdates[7] = @{
    "I" = "0.2961697031425515";
    "J" = "0.2961697031425515";
    "O" = "0.18418156894094";
    "P" = "0.18418156894094";
    "S" = "0.05454900059756759";
    "T" = "0.05454900059756759"
}

$rowUpdates[8] = @{
    "I" = "0.6932415438722715";
    "J" = "0.6932415438722715";
    "M" = "0.04154133333333333";
    "N" = "0.124624";
    "O" = "0.002517093804502335";
    "P" = "0.002517093804502335";
    "Q" = "0.1733415709724444";
    "R" = "1.560074138752";
    "S" = "0.001744953995104528";
    "T" = "0.001744953995104528"
}

$rowUpdates[9] = @{
    "I" = "0.6932415438722715";
    "J" = "0.6932415438722715";
    "O" = "0.8133013372545576";
    "P" = "0.8133013372545578";
    "S" = "0.5638142746717325";
    "T" = "0.5638142746717326"
}

$rowUpdates[10] = @{
    "I" = "0.6932415438722715";
    "J" = "0.6932415438722715";
    "O" = "0.18418156894094";
    "P" = "0.18418156894094";
    "S" = "0.1276823152054344";
    "T" = "0.1276823152054345"
}

foreach ($rowNum in $rowUpdates.Keys) {
    $cols = $rowUpdates[$rowNum]
    foreach ($col in $cols.Keys) {
        $cellRef = "$col$rowNum"
        $newValue = [double]$cols[$col]
        $ws.Range($cellRef).Value = $newValue
    }
}

Write-Host "Applied NATMI TPM update to $($rowUpdates.Keys.Count) rows"
